# Add two new "Water" characterization-factor rows for the Freshwater Use
# planetary boundary method:
#   - Water / water::fossil well            = -1e-9
#   - Water / water::ground-, long-term     = -1e-9
#
# Mirrors the commit "Add CFs for FWU PB".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows, appended right after the existing last row (row 14).
$ws.Range("A15").Value = "Water"
$ws.Range("B15").Value = "water::fossil well"
$ws.Range("C15").Value = -0.000000001

$ws.Range("A16").Value = "Water"
$ws.Range("B16").Value = "water::ground-, long-term"
$ws.Range("C16").Value = -0.000000001

# Copy the highlight formatting from the previous "new addition" row (14)
# onto the two freshly added rows so they match the rest of the sheet.
$ws.Range("A14:C14").Copy()
$ws.Range("A15:C15").PasteSpecial(-4122)
$ws.Range("A14:C14").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)

# The active selection moves down with the newly added rows.
$ws.Range("B19").Select()
